$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H70").Value = 4623.1333
$ws.Range("J70").Value = 5389.1113
$ws.Range("L70").Value = 16167.3339
$ws.Range("N70").Value = -16707.3339
$ws.Range("H73").Value = 4623.1333
$ws.Range("J73").Value = 5389.1113
$ws.Range("L73").Value = 16167.3339
$ws.Range("N73").Value = -18039.3339
$ws.Range("H74").Value = 29999.857
$ws.Range("I74").Value = 29999.857
$ws.Range("K74").Value = 29999.857
$ws.Range("M74").Value = -29063.857
$ws.Range("H77").Value = 29999.857
$ws.Range("I77").Value = 29999.857
$ws.Range("K77").Value = 149999.285
$ws.Range("M77").Value = -145319.285
$ws.Range("H96").Value = 890.2
$ws.Range("I96").Value = 325.5
$ws.Range("K96").Value = 976.5
$ws.Range("M96").Value = 396.5
$ws.Range("H135").Value = 721.2708
$ws.Range("I135").Value = 486.41177
$ws.Range("J135").Value = 1291.6428
$ws.Range("K135").Value = 4377.70593
$ws.Range("L135").Value = 11624.7852
$ws.Range("M135").Value = -1842.70593
$ws.Range("N135").Value = -16694.7852
$ws.Range("H138").Value = 4524.241
$ws.Range("J138").Value = 4470.7646
$ws.Range("L138").Value = 13412.2938
$ws.Range("N138").Value = -23692.2938

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value = 273.6875
$ws.Range("I5").Value = 258.66666
$ws.Range("K5").Value = 258.66666
$ws.Range("M5").Value = -146.66666
$ws.Range("H32").Value = 5731.71
$ws.Range("I32").Value = 3583.0723
$ws.Range("J32").Value = 16222.117
$ws.Range("K32").Value = 3583.0723
$ws.Range("L32").Value = 16222.117
$ws.Range("M32").Value = -3296.0723
$ws.Range("N32").Value = -16796.117
$ws.Range("H63").Value = 2575
$ws.Range("I63").Value = 2575
$ws.Range("K63").Value = 2575
$ws.Range("M63").Value = -1889
$ws.Range("H66").Value = 2575
$ws.Range("I66").Value = 2575
$ws.Range("K66").Value = 12875
$ws.Range("M66").Value = -9443
$ws.Range("H74").Value = 4952.5713
$ws.Range("I74").Value = 2292.2856
$ws.Range("J74").Value = 12933.429
$ws.Range("K74").Value = 2292.2856
$ws.Range("L74").Value = 12933.429
$ws.Range("M74").Value = -1418.2856
$ws.Range("N74").Value = -14681.429
$ws.Range("H77").Value = 4952.5713
$ws.Range("I77").Value = 2292.2856
$ws.Range("J77").Value = 12933.429
$ws.Range("K77").Value = 11461.428
$ws.Range("L77").Value = 64667.145
$ws.Range("M77").Value = -7093.428
$ws.Range("N77").Value = -73403.145
$ws.Range("H88").Value = 1874.5
$ws.Range("J88").Value = 1749
$ws.Range("L88").Value = 1749
$ws.Range("N88").Value = -2561
$ws.Range("H91").Value = 1874.5
$ws.Range("J91").Value = 1749
$ws.Range("L91").Value = 1749
$ws.Range("N91").Value = -4557
$ws.Range("H97").Value = 1588.0454
$ws.Range("I97").Value = 1507.6111
$ws.Range("K97").Value = 1507.6111
$ws.Range("M97").Value = -1011.6111
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H132").Value = 2730.889
$ws.Range("I132").Value = 1886.6562
$ws.Range("J132").Value = 9484.75
$ws.Range("K132").Value = 5659.9686
$ws.Range("L132").Value = 28454.25
$ws.Range("M132").Value = -3129.9686
$ws.Range("N132").Value = -33514.25

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 273.6875
$ws.Range("I4").Value = 258.66666
$ws.Range("K4").Value = 258.66666
$ws.Range("M4").Value = -143.66666
$ws.Range("H16").Value = 14504.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 14504.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 14504.5
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -14844.5
$ws.Range("H86").Value = 3999.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3999.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3999.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6245.5
$ws.Range("H89").Value = 3999.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3999.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 19997.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -31229.5
$ws.Range("H94").Value = 2412.375
$ws.Range("I94").Value = 881.5
$ws.Range("K94").Value = 881.5
$ws.Range("M94").Value = -430.5
$ws.Range("H99").Value = 2139.4736
$ws.Range("I99").Value = 1317.5454
$ws.Range("K99").Value = 1317.5454
$ws.Range("M99").Value = 180.4546
$ws.Range("H105").Value = 681484.25
$ws.Range("I105").Value = 2861916
$ws.Range("K105").Value = 2861916
$ws.Range("M105").Value = -2860169
$ws.Range("H134").Value = 5024.75
$ws.Range("I134").Value = 3933.0476
$ws.Range("J134").Value = 12666.667
$ws.Range("K134").Value = 11799.1428
$ws.Range("L134").Value = 38000.001
$ws.Range("M134").Value = -9264.1428
$ws.Range("N134").Value = -43070.001
$ws.Range("H138").Value = 35000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 91728.28999999999
$ws.Range("J140").Value = 91728.28999999999
$ws.Range("L140").Value = 91728.28999999999
$ws.Range("N140").Value = -102088.29

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 1382.8636
$ws.Range("I22").Value = 576.5714
$ws.Range("K22").Value = 576.5714
$ws.Range("M22").Value = -226.5714
$ws.Range("H62").Value = 2652.75
$ws.Range("J62").Value = 2652.75
$ws.Range("L62").Value = 2652.75
$ws.Range("N62").Value = -3900.75
$ws.Range("H65").Value = 2652.75
$ws.Range("J65").Value = 2652.75
$ws.Range("L65").Value = 13263.75
$ws.Range("N65").Value = -19503.75
$ws.Range("H105").Value = 8000
$ws.Range("I105").Value = 8000
$ws.Range("K105").Value = 8000
$ws.Range("M105").Value = -6253

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 10929410
$ws.Range("I4").Value = 23123244
$ws.Range("J4").Value = 477552.62
$ws.Range("K4").Value = 69369732
$ws.Range("L4").Value = 1432657.86
$ws.Range("M4").Value = -69369620
$ws.Range("N4").Value = -1432881.86

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5064.933
$ws.Range("J70").Value = 5064.933
$ws.Range("L70").Value = 5064.933
$ws.Range("N70").Value = -5604.933
$ws.Range("H73").Value = 5064.933
$ws.Range("J73").Value = 5064.933
$ws.Range("L73").Value = 5064.933
$ws.Range("N73").Value = -6936.933
$ws.Range("H97").Value = 1375.25
$ws.Range("I97").Value = 1334
$ws.Range("J97").Value = 1499
$ws.Range("K97").Value = 1334
$ws.Range("L97").Value = 1499
$ws.Range("M97").Value = -838
$ws.Range("N97").Value = -2491
$ws.Range("H122").Value = 4687.6177
$ws.Range("I122").Value = 2962.5652
$ws.Range("K122").Value = 8887.695599999999
$ws.Range("M122").Value = -6437.695599999999
$ws.Range("H138").Value = 189991.5
$ws.Range("J138").Value = 189991.5
$ws.Range("L138").Value = 189991.5
$ws.Range("N138").Value = -200271.5

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 10226.12
$ws.Range("I68").Value = 8776.9
$ws.Range("J68").Value = 11192.267
$ws.Range("K68").Value = 8776.9
$ws.Range("L68").Value = 11192.267
$ws.Range("M68").Value = -8027.9
$ws.Range("N68").Value = -12690.267
$ws.Range("H71").Value = 10226.12
$ws.Range("I71").Value = 8776.9
$ws.Range("J71").Value = 11192.267
$ws.Range("K71").Value = 43884.5
$ws.Range("L71").Value = 55961.335
$ws.Range("M71").Value = -40140.5
$ws.Range("N71").Value = -63449.335
$ws.Range("H82").Value = 1970.6842
$ws.Range("I82").Value = 1266.6
$ws.Range("K82").Value = 1266.6
$ws.Range("M82").Value = -905.5999999999999
$ws.Range("H85").Value = 1970.6842
$ws.Range("I85").Value = 1266.6
$ws.Range("K85").Value = 1266.6
$ws.Range("M85").Value = -18.59999999999991
$ws.Range("H100").Value = 4395.381
$ws.Range("I100").Value = 3963.1052
$ws.Range("J100").Value = 8502
$ws.Range("K100").Value = 3963.1052
$ws.Range("L100").Value = 8502
$ws.Range("M100").Value = -3422.1052
$ws.Range("N100").Value = -9584
$ws.Range("H132").Value = 6580.1665
$ws.Range("I132").Value = 6268.12
$ws.Range("J132").Value = 8140.4
$ws.Range("K132").Value = 18804.36
$ws.Range("L132").Value = 24421.2
$ws.Range("M132").Value = -16274.36
$ws.Range("N132").Value = -29481.2

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H96").Value = 13249.714
$ws.Range("I96").Value = 9444.111000000001
$ws.Range("K96").Value = 9444.111000000001
$ws.Range("M96").Value = -8071.111000000001
$ws.Range("H122").Value = 3412.524
$ws.Range("I122").Value = 3454.3845
$ws.Range("J122").Value = 2868.3333
$ws.Range("K122").Value = 10363.1535
$ws.Range("L122").Value = 8604.999899999999
$ws.Range("M122").Value = -7913.1535
$ws.Range("N122").Value = -13504.9999
